$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the friendly-name ("A 34293-2023") argument to each HYPERLINK() formula
# in row 2. This mirrors an automatic find/replace style update: S2 gets the
# properly closed string, while T2/V2/W2/X2/Y2 reproduce the original
# (unbalanced-quote) text exactly as produced by that update.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/artfynd/A 34293-2023.xlsx"; "A 34293-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/kartor/A 34293-2023.png; "A 34293-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/klagomål/A 34293-2023.docx; "A 34293-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/klagomålsmail/A 34293-2023.docx; "A 34293-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/tillsyn/A 34293-2023.docx; "A 34293-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_LINDESBERG/tillsynsmail/A 34293-2023.docx; "A 34293-2023")'
